$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.808.65'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '3.531.40'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '''621.34'
$ws.Range("E5").Value = '  +4.38%  '
$ws.Range("D6").Value = '''172.67'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").Value = '3.528.57'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.608'
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +2.41%  '
$ws.Range("D11").Value = '''7.25'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("D12").Value = '''0.585'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '''46.23'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '4.102.58'
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '''8.45'
$ws.Range("E16").Value = '  +2.01%  '
$ws.Range("D17").Value = '''607.10'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '3.533.71'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '70.906.98'
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").Value = '''17.66'
$ws.Range("E21").Value = '  +2.53%  '
$ws.Range("D22").Value = '''0.881'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").Value = '''15.70'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("D25").Value = '''97.82'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").Value = '''33.67'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '''9.11'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '''6.84'
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("D35").Value = '''616.86'
$ws.Range("E35").Value = '  -7.71%  '
$ws.Range("D36").Value = '''0.0503'
$ws.Range("E36").Value = '  +6.56%  '
$ws.Range("D37").Value = '''10.89'
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").Value = '''0.0996'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '''56.96'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '''3.39'
$ws.Range("E41").Value = '  -4.79%  '
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = '3.345.97'
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D44").Value = '0.0₃0729'
$ws.Range("E44").Value = '  +5.25%  '
$ws.Range("D45").Value = '''0.311'
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("D46").Value = '''2.88'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").Value = '''31.79'
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").Value = '''134.05'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("B51").Value = 'USDe'
$ws.Range("C51").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  -0.01%  '
